# Fill in the sample/profile row (row 2) of the profile template and
# attach data validation rules to each of its cells, matching the
# "Male" example row shown in the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 sample data -----------------------------------------------
# A2: Phone Number
$ws.Range("A2").Value2 = 9381153612
# B2: Date of Birth (10-Jun-2003), keep existing date-format style (s="2")
$ws.Range("B2").Value2 = 37782
# C2: Gender
$ws.Range("C2").Value = "Male"
# E2: Income Range
$ws.Range("E2").Value2 = 344

# --- Data validations ---------------------------------------------------
# C2 - Gender list
$ws.Range("C2").Validation.Add(3, 1, 1, '"Male,Femal,Other"')
# D2 - Marital Status list
$ws.Range("D2").Validation.Add(3, 1, 1, '"Single,Married,Divorced,Widowed"')
# E2 - Income Range whole number between 1 and 1E+46
$ws.Range("E2").Validation.Add(1, 1, 1, "1", "1E+46")
# A2 - Phone Number custom rule (numeric, 10 digits)
$ws.Range("A2").Validation.Add(7, 1, 1, "AND(ISNUMBER(A2),LEN(A2)=10)")
# B2 - Date of Birth between 1 and TODAY()
$ws.Range("B2").Validation.Add(4, 1, 1, "1", "TODAY()")

# --- Selection moves to F4 ---------------------------------------------
$ws.Range("F4").Select()
